$wb = $excel.ActiveWorkbook

# Both the "展览" (sheet 1) and "全部类型" (sheet 4) worksheets hold the same
# event rows and need the same updated numbers (refreshed scrape output).
$sheetIndexes = @(1, 4)

foreach ($idx in $sheetIndexes) {
    $ws = $wb.Worksheets.Item($idx)

    # Row 2: interest count (F) and min price text (G)
    $ws.Range("F2").Value = 1977
    $g2 = $ws.Range("G2")
    $g2.Value = "'70"
    $g2.Style = "Normal"

    # Row 4: interest count (F)
    $ws.Range("F4").Value = 269

    # Row 5: interest count (F)
    $ws.Range("F5").Value = 10487

    # Row 6: interest count (F)
    $ws.Range("F6").Value = 9126

    # Row 8: interest count (F)
    $ws.Range("F8").Value = 664

    # Row 12: interest count (F)
    $ws.Range("F12").Value = 9400

    # Row 13: interest count (F)
    $ws.Range("F13").Value = 2412

    # Row 15: interest count (F)
    $ws.Range("F15").Value = 57

    # Row 16: interest count (F)
    $ws.Range("F16").Value = 356

    # Row 17: interest count (F)
    $ws.Range("F17").Value = 10613

    # Row 18: interest count (F)
    $ws.Range("F18").Value = 10675

    # Row 20: interest count (F)
    $ws.Range("F20").Value = 5
}
